$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: update "Have" quantity, recalculated "Need", and add a note in the
#     new Buy column. (Set this text first so the shared-string table ends up
#     in the same order the original author typed things in.) ---
$ws.Range("E25").Value = 12
$ws.Rows.Item(25).RowHeight = 28.8

$g25 = $ws.Range("G25")
$g25.NumberFormat = "General"
$g25.WrapText = $true
$g25.Value = "200mA hold, 400mA trip"

# --- Header row: G1 "Buy" label, K1 label text change ---
$ws.Range("G1").NumberFormat = "General"
$ws.Range("G1").WrapText = $true
$ws.Range("G1").Value = "Buy"

$ws.Range("K1").Value = "Needed per board"

# --- New "Buy" quantities in column G for various rows ---
$buyCells = @{
    "G3"  = 2
    "G4"  = 3
    "G5"  = 3
    "G7"  = 15
    "G17" = 5
    "G18" = 5
    "G20" = 6
    "G21" = 3
    "G22" = 3
    "G23" = 3
    "G24" = 2
}

foreach ($addr in $buyCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "General"
    $cell.WrapText = $true
    $cell.Value = $buyCells[$addr]
}

# --- Column G width (widen to fit the "Buy" column) ---
$ws.Columns.Item(7).ColumnWidth = 13.6

# --- Rows 20-22: column F now uses the same shared IF() formula as the rest ---
$ws.Range("F20").Formula = "=IF(D20>E20, D20-E20, 0)"
$ws.Range("F21").Formula = "=IF(D21>E21, D21-E21, 0)"
$ws.Range("F22").Formula = "=IF(D22>E22, D22-E22, 0)"

# --- View state: scroll back to top-left and move the active selection ---
$ws.Range("H12").Select() | Out-Null
